$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 8335210
$ws.Range("J19").Value = 3221.4285
$ws.Range("L19").Value = 3221.4285
$ws.Range("N19").Value = -3571.4285

# Row 64
$ws.Range("H64").Value = 3068.4443
$ws.Range("I64").Value = 2984.138
$ws.Range("J64").Value = 3221.25
$ws.Range("K64").Value = 2984.138
$ws.Range("L64").Value = 3221.25
$ws.Range("M64").Value = -2736.138
$ws.Range("N64").Value = -3717.25

# Row 67
$ws.Range("H67").Value = 3068.4443
$ws.Range("I67").Value = 2984.138
$ws.Range("J67").Value = 3221.25
$ws.Range("K67").Value = 2984.138
$ws.Range("L67").Value = 3221.25
$ws.Range("M67").Value = -2126.138
$ws.Range("N67").Value = -4937.25

# Row 70
$ws.Range("H70").Value = 2340.7407
$ws.Range("I70").Value = 952.3333
$ws.Range("J70").Value = 2737.4285
$ws.Range("K70").Value = 2856.9999
$ws.Range("L70").Value = 8212.2855
$ws.Range("M70").Value = -2586.9999
$ws.Range("N70").Value = -8752.2855

# Row 73
$ws.Range("H73").Value = 2340.7407
$ws.Range("I73").Value = 952.3333
$ws.Range("J73").Value = 2737.4285
$ws.Range("K73").Value = 2856.9999
$ws.Range("L73").Value = 8212.2855
$ws.Range("M73").Value = -1920.9999
$ws.Range("N73").Value = -10084.2855

# Row 103
$ws.Range("H103").Value = 40001932
$ws.Range("I103").Value = 125000616
$ws.Range("J103").Value = 2553.1765
$ws.Range("K103").Value = 375001848
$ws.Range("L103").Value = 7659.529500000001
$ws.Range("M103").Value = -375001262
$ws.Range("N103").Value = -8831.529500000001

# Row 116
$ws.Range("H116").Value = 9378914
$ws.Range("I116").Value = 4446754.5
$ws.Range("J116").Value = 16777153
$ws.Range("K116").Value = 4446754.5
$ws.Range("L116").Value = 16777153
$ws.Range("M116").Value = -4443312.5
$ws.Range("N116").Value = -16784037

# Row 127
$ws.Range("H127").Value = 1748.8667
$ws.Range("I127").Value = 533.3
$ws.Range("J127").Value = 4180
$ws.Range("K127").Value = 1599.9
$ws.Range("L127").Value = 12540
$ws.Range("M127").Value = 3360.1
$ws.Range("N127").Value = -22460

# Row 132
$ws.Range("H132").Value = 2671334.8
$ws.Range("I132").Value = 597921.0600000001
$ws.Range("J132").Value = 12347266
$ws.Range("K132").Value = 1793763.18
$ws.Range("L132").Value = 37041798
$ws.Range("M132").Value = -1791233.18
$ws.Range("N132").Value = -37046858

# Row 134
$ws.Range("H134").Value = 142711.58
$ws.Range("J134").Value = 142711.58
$ws.Range("L134").Value = 142711.58
$ws.Range("N134").Value = -152851.58

$ws = $wb.Worksheets.Item("ARM")
# Row 37
$ws.Range("H37").Value = 11982.692
$ws.Range("I37").Value = 6489
$ws.Range("J37").Value = 16691.572
$ws.Range("K37").Value = 6489
$ws.Range("L37").Value = 16691.572
$ws.Range("M37").Value = -6216
$ws.Range("N37").Value = -17237.572

# Row 74
$ws.Range("H74").Value = 52875440
$ws.Range("I74").Value = 45455324
$ws.Range("J74").Value = 76195810
$ws.Range("K74").Value = 45455324
$ws.Range("L74").Value = 76195810
$ws.Range("M74").Value = -45454450
$ws.Range("N74").Value = -76197558

# Row 77
$ws.Range("H77").Value = 52875440
$ws.Range("I77").Value = 45455324
$ws.Range("J77").Value = 76195810
$ws.Range("K77").Value = 227276620
$ws.Range("L77").Value = 380979050
$ws.Range("M77").Value = -227272252
$ws.Range("N77").Value = -380987786

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1492560.4
$ws.Range("I31").Value = 2689738.8
$ws.Range("J31").Value = 8059.28
$ws.Range("K31").Value = 2689738.8
$ws.Range("L31").Value = 8059.28
$ws.Range("M31").Value = -2689443.8
$ws.Range("N31").Value = -8649.279999999999

# Row 34
$ws.Range("H34").Value = 1492560.4
$ws.Range("I34").Value = 2689738.8
$ws.Range("J34").Value = 8059.28
$ws.Range("K34").Value = 2689738.8
$ws.Range("L34").Value = 8059.28
$ws.Range("M34").Value = -2689536.8
$ws.Range("N34").Value = -8463.279999999999

# Row 51
$ws.Range("H51").Value = 9400.4
$ws.Range("J51").Value = 9400.4
$ws.Range("L51").Value = 9400.4
$ws.Range("N51").Value = -10872.4

# Row 59
$ws.Range("H59").Value = 16700.334
$ws.Range("J59").Value = 16700.334
$ws.Range("L59").Value = 16700.334
$ws.Range("N59").Value = -18990.334

# Row 60
$ws.Range("H60").Value = 6160.4
$ws.Range("J60").Value = 7200.5
$ws.Range("L60").Value = 7200.5
$ws.Range("N60").Value = -8222.5

# Row 61
$ws.Range("H61").Value = 9400.4
$ws.Range("J61").Value = 9400.4
$ws.Range("L61").Value = 9400.4
$ws.Range("N61").Value = -10096.4

# Row 68
$ws.Range("H68").Value = 17999.727
$ws.Range("I68").Value = 20666.334
$ws.Range("J68").Value = 16999.75
$ws.Range("K68").Value = 20666.334
$ws.Range("L68").Value = 16999.75
$ws.Range("M68").Value = -19917.334
$ws.Range("N68").Value = -18497.75

# Row 71
$ws.Range("H71").Value = 17999.727
$ws.Range("I71").Value = 20666.334
$ws.Range("J71").Value = 16999.75
$ws.Range("K71").Value = 61999.00199999999
$ws.Range("L71").Value = 50999.25
$ws.Range("M71").Value = -58255.00199999999
$ws.Range("N71").Value = -58487.25

# Row 74
$ws.Range("H74").Value = 16283.833
$ws.Range("J74").Value = 16283.833
$ws.Range("L74").Value = 16283.833
$ws.Range("N74").Value = -18031.833

# Row 77
$ws.Range("H77").Value = 16283.833
$ws.Range("J77").Value = 16283.833
$ws.Range("L77").Value = 48851.499
$ws.Range("N77").Value = -57587.499

# Row 135
$ws.Range("H135").Value = 49966.668
$ws.Range("J135").Value = 49966.668
$ws.Range("L135").Value = 49966.668
$ws.Range("N135").Value = -60106.668

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Range("H3").Value = 8285.714
$ws.Range("I3").Value = 8000
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 24000
$ws.Range("L3").Value = 30000
$ws.Range("M3").Value = -23888
$ws.Range("N3").Value = -30224

$ws = $wb.Worksheets.Item("LTW")
# Row 82
$ws.Range("H82").Value = 4682.2104
$ws.Range("I82").Value = 1172.375
$ws.Range("J82").Value = 7234.8184
$ws.Range("K82").Value = 1172.375
$ws.Range("L82").Value = 7234.8184
$ws.Range("M82").Value = -811.375
$ws.Range("N82").Value = -7956.8184

# Row 85
$ws.Range("H85").Value = 4682.2104
$ws.Range("I85").Value = 1172.375
$ws.Range("J85").Value = 7234.8184
$ws.Range("K85").Value = 1172.375
$ws.Range("L85").Value = 7234.8184
$ws.Range("M85").Value = 75.625
$ws.Range("N85").Value = -9730.8184

# Row 132
$ws.Range("H132").Value = 8405521
$ws.Range("I132").Value = 12988251
$ws.Range("J132").Value = 3850.6667
$ws.Range("K132").Value = 38964753
$ws.Range("L132").Value = 11552.0001
$ws.Range("M132").Value = -38962223
$ws.Range("N132").Value = -16612.0001

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 39914.5
$ws.Range("J46").Value = 39914.5
$ws.Range("L46").Value = 39914.5
$ws.Range("N46").Value = -40376.5

# Row 134
$ws.Range("H134").Value = 39914.5
$ws.Range("J134").Value = 39914.5
$ws.Range("L134").Value = 119743.5
$ws.Range("N134").Value = -124813.5

# Row 136
$ws.Range("H136").Value = 1664.641
$ws.Range("I136").Value = 1300.1428
$ws.Range("J136").Value = 1744.375
$ws.Range("K136").Value = 3900.4284
$ws.Range("L136").Value = 5233.125
$ws.Range("M136").Value = -1350.4284
$ws.Range("N136").Value = -10333.125
